$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ran the averaged-intensity notebook including three new spiral sampling
# schemes (Spiral-90deg-10rot-5space / -15rot-5space / -10rot-3space). The table
# of schemes got re-sorted, so rows 10-16 now hold different scheme rows and the
# sheet grows from 16 to 19 rows (A1:P16 -> A1:P19).

# New row/label/data mapping for rows 10-19
$rowData = @{
    10 = @{ A = 8; B = "Gaussian-Quadrature"; Vals = @(1.010501375714796, 0.9862014279361111, 0.9964998363637791, 0.9879085223676277, 1.010501375714796, 0.9862014279361111, 0.995867728876934, 0.9923383171122281, 0.9982352941176471, 0.9805882352941176, 1.010501375714796, 0.9913506321499451, 0.9952777905955785, 0.9935175922229051) }
    11 = @{ A = 9; B = "Spiral-90deg-10rot-5space"; Vals = @(1.009141297048027, 0.962236934594794, 1.003775096615783, 0.9863061149377539, 1.009141297048027, 0.962236934594794, 1.004338935143897, 0.9863787482844377, 1.00305669878758, 0.9749059388749645, 1.009141297048027, 0.9830060156052882, 0.9903648607990895, 0.9912674705359046) }
    12 = @{ A = 10; B = "Spiral-90deg-15rot-5space"; Vals = @(1.008969542755898, 0.9623477900679072, 1.003724105759154, 0.9863474077727332, 1.008969542755898, 0.9623477900679072, 1.004295705195686, 0.986396280864007, 1.00303111124825, 0.9750183032758598, 1.008969542755898, 0.9830359479135306, 0.9903472115889231, 0.9912662808674368) }
    13 = @{ A = 11; B = "Spiral-90deg-10rot-3space"; Vals = @(1.009123139392663, 0.9622327267136841, 1.003756588104599, 0.9863395991945588, 1.009123139392663, 0.9622327267136841, 1.004328664447375, 0.9863840726755867, 1.003044740783079, 0.9749474487110216, 1.009123139392663, 0.9829946574091415, 0.9903630133513761, 0.9912696225028208) }
    14 = @{ A = 12; B = "NoRotation-tilt60deg"; Vals = @(1.000415999999999, 0.9612399999999994, 1.011467999999999, 0.9801960000000012, 1.000415999999999, 0.9612399999999994, 1.009592, 0.9835400000000006, 1.005628, 0.9660279999999996, 1.000415999999999, 0.9863539999999992, 0.9883299999999995, 0.9897634999999997) }
    15 = @{ A = 13; B = "Rotation-NoTilt"; Vals = @(1, 0.9438874999999985, 1.02, 0.97, 1, 0.9438874999999985, 1.02, 0.98, 1.01, 0.9520750000000017, 1, 0.9819437499999992, 0.9834718749999996, 0.9869953124999999) }
    16 = @{ A = 14; B = "Rotation-60detTilt"; Vals = @(0.9973799755776024, 0.9640421148671958, 1.008466805555199, 0.9801391959040007, 0.9973799755776024, 0.9640421148671958, 1.008646414131197, 0.9852232388608003, 1.002462425088002, 0.9696016136192026, 0.9973799755776024, 0.9862544602111973, 0.9875070229759995, 0.9894952229503999) }
    17 = @{ A = 15; B = "HexGrid-90degTilt5degRes"; Vals = @(0.9926564473382259, 0.9925726426016895, 0.9919755244143033, 0.992528936875294, 0.9926564473382259, 0.9925726426016895, 0.9926475533560726, 0.9923321043988091, 0.9926535869352376, 0.9923073296038305, 0.9926564473382259, 0.9922740835079964, 0.9924333878073781, 0.9924592656904327) }
    18 = @{ A = 16; B = "HexGrid-90degTilt22p5degRes"; Vals = @(0.9921897788832952, 0.9928849056754979, 0.9905216880934852, 0.9925158965746071, 0.9921897788832952, 0.9928849056754979, 0.9939433958500613, 0.9930541355605793, 0.9910340834547853, 0.99345164867309, 0.9921897788832952, 0.9917032968844915, 0.9920280673067213, 0.9924494415956752) }
    19 = @{ A = 17; B = "HexGrid-60degTilt5degRes"; Vals = @(0.989184923041861, 0.9986449157464137, 0.9908984188236678, 0.9937226641067702, 0.989184923041861, 0.9986449157464137, 0.9906355018040066, 0.9936083723699995, 0.9909715415884404, 0.9965045309483396, 0.989184923041861, 0.9947716672850407, 0.9931127304296783, 0.9930213585536873) }
}

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P")

# A3 already carries the bold/centered/bordered "index" style used down column A;
# copy it onto the three brand-new rows (17-19) before the A1:A16 values get
# rewritten below, then overwrite every cell with the post-rerun values.
$ws.Range("A3").Copy($ws.Range("A17"))
$ws.Range("A3").Copy($ws.Range("A18"))
$ws.Range("A3").Copy($ws.Range("A19"))

foreach ($r in ($rowData.Keys | Sort-Object)) {
    $info = $rowData[$r]
    $ws.Range("A$r").Value = $info.A
    $ws.Range("B$r").Value = $info.B
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $info.Vals[$i]
    }
}